$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.040.69"
$ws.Range("E2").Value = "  +0.48%  "

$ws.Range("D3").Value = "1.677.08"
$ws.Range("E3").Value = "  +0.52%  "

$ws.Range("E4").Value = "  +0.08%  "

$ws.Range("D5").Value = "'215.47"
$ws.Range("E5").Value = "  +0.00%  "

$ws.Range("E6").Value = "  -0.97%  "

$ws.Range("E7").Value = "  +0.05%  "

$ws.Range("D8").Value = "'0.254"
$ws.Range("E8").Value = "  +1.88%  "

$ws.Range("D9").Value = "'21.34"
$ws.Range("E9").Value = "  +5.38%  "

$ws.Range("E10").Value = "  +0.32%  "

$ws.Range("D11").Value = "'0.0888"
$ws.Range("E11").Value = "  -0.70%  "

$ws.Range("D12").Value = "1.911.40"
$ws.Range("E12").Value = "  +0.39%  "

$ws.Range("D13").Value = "1.697.46"
$ws.Range("E13").Value = "  +1.25%  "

$ws.Range("E14").Value = "  +0.82%  "

$ws.Range("D15").Value = "'0.535"
$ws.Range("E15").Value = "  +1.70%  "

$ws.Range("D16").Value = "'66.32"

$ws.Range("D17").Value = "27.036.24"
$ws.Range("E17").Value = "  +0.41%  "

$ws.Range("D18").Value = "'8.14"
$ws.Range("E18").Value = "  +1.80%  "

$ws.Range("D19").Value = "'235.18"
$ws.Range("E19").Value = "  +0.59%  "

$ws.Range("D20").Value = "0.0₃0735"
$ws.Range("E20").Value = "  +0.40%  "

$ws.Range("E21").Value = "  +0.03%  "

$ws.Range("E22").Value = "  +1.14%  "

$ws.Range("D23").Value = "'9.25"
$ws.Range("E23").Value = "  +1.31%  "

$ws.Range("D24").Value = "'2.11"
$ws.Range("E24").Value = "  -3.48%  "

$ws.Range("D25").Value = "'147.33"
$ws.Range("E25").Value = "  +0.79%  "

$ws.Range("E26").Value = "  +1.79%  "

$ws.Range("D27").Value = "'16.46"
$ws.Range("E27").Value = "  +3.62%  "

$ws.Range("E28").Value = "  -1.45%  "

$ws.Range("E29").Value = "  +0.18%  "

$ws.Range("E30").Value = "  +0.19%  "

$ws.Range("E31").Value = "  +0.01%  "

$ws.Range("E32").Value = "  +0.59%  "

$ws.Range("D33").Value = "1.539.80"
$ws.Range("E33").Value = "  +5.93%  "

$ws.Range("D34").Value = "'3.19"
$ws.Range("E34").Value = "  +2.07%  "

$ws.Range("E35").Value = "  +4.57%  "

$ws.Range("E36").Value = "  -1.13%  "

$ws.Range("E37").Value = "  +0.01%  "

$ws.Range("E38").Value = "  +1.00%  "

$ws.Range("D39").Value = "'0.0174"
$ws.Range("E39").Value = "  +2.14%  "

$ws.Range("E40").Value = "  +6.78%  "

$ws.Range("E41").Value = "  +0.03%  "

$ws.Range("D42").Value = "'67.82"
$ws.Range("E42").Value = "  +2.64%  "

$ws.Range("E43").Value = "  -3.69%  "

$ws.Range("E44").Value = "  -1.17%  "

$ws.Range("D45").Value = "1.818.31"
$ws.Range("E45").Value = "  -0.10%  "

$ws.Range("E46").Value = "  -0.62%  "

$ws.Range("D47").Value = "'90.32"
$ws.Range("E47").Value = "  -0.42%  "

$ws.Range("E48").Value = "  +2.79%  "

$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("E50").Value = "  +1.69%  "

$ws.Range("E51").Value = "  +5.47%  "
